# Update res_bus/vm_pu.xlsx values for the 380 kV case (Case_3_64)
# Rows 2-25 (bus indices 0-23), columns B,C,D,E,F,I,J,K,L,M,N
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
    2 = @(1.02, 1.035267931653078, 1.044069565816396, 1.034341655449683, 1.05282291906517, 1.038690198370736, 1.040382636490932, 1.046841652950761, 1.037141406624156, 1.055570587335104, 1.041860099168306)
    3 = @(1.02, 1.036141159482062, 1.044788307787412, 1.035082130898107, 1.053773996054072, 1.038901160308511, 1.040899726929264, 1.047371799504011, 1.03769123293269, 1.056334232298374, 1.042377923934371)
    4 = @(1.02, 1.036706773485832, 1.045253893726203, 1.035562136944138, 1.054390455769006, 1.039036762942801, 1.041234247680258, 1.047714679348877, 1.038047209457395, 1.056828785558063, 1.042712919743201)
    5 = @(1.02, 1.036944694408423, 1.045449746989361, 1.035764137870276, 1.054649864993586, 1.039093553267368, 1.041374861988907, 1.047858786426942, 1.038196909143259, 1.057036795775987, 1.042853733740298)
    6 = @(1.02, 1.036984650359682, 1.045482638659276, 1.035798066769526, 1.054693435530414, 1.039103075867855, 1.041398470648914, 1.047882980281842, 1.038222047108451, 1.057071727430246, 1.04287737592731)
    7 = @(1.02, 1.036709952061061, 1.045256510253904, 1.035564835281505, 1.054393921029216, 1.039037522631445, 1.041236126648526, 1.047716605071846, 1.038049209568188, 1.056831564609556, 1.04271480137982)
    8 = @(1.02, 1.035562922735283, 1.044312361094232, 1.034591721695937, 1.053144121730346, 1.038761680707023, 1.040557403713378, 1.047020850954146, 1.037327180401601, 1.055828576077271, 1.042035114580254)
    9 = @(1.02, 1.033546192732641, 1.042652640191662, 1.032883698407635, 1.050949924949621, 1.03826871798316, 1.039360911061332, 1.045793667838288, 1.036056475586248, 1.054064493177221, 1.040836922771402)
    10 = @(1.02, 1.03220480204414, 1.041548937641578, 1.031749636611004, 1.049492677133147, 1.037935479170884, 1.038562983266635, 1.044974826966272, 1.035210490286508, 1.052890750636744, 1.040037861827713)
    11 = @(1.02, 1.031624716588843, 1.041071701166554, 1.031259691757649, 1.048863009647798, 1.037790100733798, 1.038217422792985, 1.044620104352176, 1.034844457842822, 1.052383073860365, 1.039691810618554)
    12 = @(1.02, 1.031409360397334, 1.040894537161199, 1.031077873056203, 1.048629324633135, 1.037735938500183, 1.038089059491466, 1.044488321869361, 1.034708541119317, 1.052194585907447, 1.039563265026423)
    13 = @(1.02, 1.031455549879085, 1.04093253476083, 1.031116866127297, 1.048679441716194, 1.037747563813404, 1.03811659414296, 1.04451659066774, 1.034737693705819, 1.052235013306989, 1.03959083878028)
    14 = @(1.02, 1.031606912858164, 1.041057054626782, 1.031244659117442, 1.048843689051022, 1.037785626972067, 1.03820681237372, 1.044609211627396, 1.034833222018818, 1.052367491628571, 1.039681185131277)
    15 = @(1.02, 1.031700187649287, 1.041133789105288, 1.031323418991322, 1.048944914008103, 1.037809057464804, 1.038262397941861, 1.044666275484006, 1.03489208603594, 1.052449127308361, 1.0397368496373)
    16 = @(1.02, 1.032243316229705, 1.041580624592582, 1.031782176194578, 1.049534494261263, 1.037945104677105, 1.038585915967899, 1.04499836548756, 1.035234788759958, 1.052924455425934, 1.040060827096044)
    17 = @(1.02, 1.032584207007966, 1.041861094275714, 1.03207024118914, 1.049904679731415, 1.038030153766412, 1.038788837025725, 1.045206635017147, 1.035449834241666, 1.053222767544564, 1.040264036325046)
    18 = @(1.02, 1.032783114676087, 1.042024752421331, 1.03223837174848, 1.050120730910531, 1.038079656845172, 1.038907192306998, 1.045328099677791, 1.035575294026248, 1.053396822034577, 1.040382559684395)
    19 = @(1.02, 1.032850949195391, 1.042080566583113, 1.032295718043562, 1.05019442048411, 1.038096518352085, 1.038947547493446, 1.04536951333422, 1.035618077179306, 1.053456179246277, 1.040422972179836)
    20 = @(1.02, 1.032547625212552, 1.041830995805987, 1.032039323449373, 1.049864949055739, 1.038021039620455, 1.038767066040009, 1.045184291255273, 1.03542675905427, 1.0531907558734, 1.040242234422034)
    21 = @(1.02, 1.031562337092783, 1.041020383777197, 1.031207022598622, 1.048795316728379, 1.037774422793279, 1.038180245530812, 1.044581937680594, 1.034805090086395, 1.052328477684482, 1.039654580560404)
    22 = @(1.02, 1.030943503975865, 1.040511315793982, 1.030684698492198, 1.048123963209446, 1.037618426824902, 1.037811249659828, 1.044203083280972, 1.034414477453636, 1.051786826330372, 1.039285060672959)
    23 = @(1.02, 1.031271496422807, 1.04078112533564, 1.030961499271886, 1.048479749265826, 1.037701211958823, 1.038006864627486, 1.044403933165247, 1.034621523985354, 1.052073918351415, 1.039480953436308)
    24 = @(1.02, 1.032564154729411, 1.041844595808471, 1.032053293514352, 1.049882901241782, 1.038025158234015, 1.038776903427339, 1.045194387487307, 1.0354371856537, 1.053205220402743, 1.040252085779582)
    25 = @(1.02, 1.034067026139857, 1.043081235135971, 1.033324456122477, 1.051516206353915, 1.038396973717456, 1.039670284930376, 1.046111056023803, 1.036384785714688, 1.054520149605759, 1.041146735986826)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}

Write-Output "Updated $($data.Keys.Count) rows"
